# Edit the "債務" (debt) sheet (sheet7) to:
#  - turn row 1 (which used to just duplicate the first data row) into the
#    real column-header row, adding two new header labels: "species" and
#    "debtor"
#  - append the standard trailing columns (property_category / category /
#    date / legislator_name / legislator_id / source_file / index) that are
#    already present on every other sheet in this workbook, filling them in
#    with the same values used elsewhere ("debt" / "normal" / "2012-02-01" /
#    "吳宜臻" / 1735 / "tmp2691" / <row index>)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(7)

# ---- Header row (row 1) ----
$ws.Cells.Item(1,2).Value = "species"
$ws.Cells.Item(1,3).Value = "debtor"
$ws.Cells.Item(1,4).Value = "owner"
$ws.Cells.Item(1,5).Value = "total"
$ws.Cells.Item(1,6).Value = "register_date"
$ws.Cells.Item(1,7).Value = "register_reason"
$ws.Cells.Item(1,8).Value = "property_category"
$ws.Cells.Item(1,9).Value = "category"
$ws.Cells.Item(1,10).Value = "date"
$ws.Cells.Item(1,11).Value = "legislator_name"
$ws.Cells.Item(1,12).Value = "legislator_id"
$ws.Cells.Item(1,13).Value = "source_file"
$ws.Cells.Item(1,14).Value = "index"

# New trailing header cells (H1:N1) need the bold/boxed header style; copy it
# over from one of the pre-existing header cells (B1).
$ws.Cells.Item(1,2).Copy() | Out-Null
$ws.Range($ws.Cells.Item(1,8), $ws.Cells.Item(1,14)).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---- Data rows (rows 2-4): columns A-G are unchanged, only H:N get added ----
for ($r = 2; $r -le 4; $r++) {
    $ws.Cells.Item($r,8).Value = "debt"
    $ws.Cells.Item($r,9).Value = "normal"
    $ws.Cells.Item($r,11).Value = "吳宜臻"
    $ws.Cells.Item($r,12).Value = 1735
    $ws.Cells.Item($r,13).Value = "tmp2691"
    $ws.Cells.Item($r,14).Value = 100 + $r
}

# The "date" column (J) holds the literal text "2012-02-01"; assigning that
# string directly makes Excel auto-convert it into a date serial number, so
# instead enter it as a formula returning the text and then convert the
# formula results down to plain values (this keeps the cell's number format
# untouched, i.e. no new style gets introduced).
$dateRng = $ws.Range($ws.Cells.Item(2,10), $ws.Cells.Item(4,10))
$dateRng.Formula = '="2012-02-01"'
$dateRng.Copy() | Out-Null
$dateRng.PasteSpecial(-4163) | Out-Null   # xlPasteValues
$excel.CutCopyMode = 0

# Copy the existing data-row style (from B2, which already carries it) onto
# the newly added trailing cells (H2:N4) so they are formatted like the rest
# of the table instead of picking up the default style.
$ws.Cells.Item(2,2).Copy() | Out-Null
$ws.Range($ws.Cells.Item(2,8), $ws.Cells.Item(4,14)).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
